$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 137.58333
$ws.Range("I6").Value = 131
$ws.Range("J6").Value = 146.8
$ws.Range("K6").Value = 393
$ws.Range("L6").Value = 440.4
$ws.Range("M6").Value = -281
$ws.Range("N6").Value = -664.4000000000001
$ws.Range("H31").Value = 3150.5
$ws.Range("I31").Value = 3150.5
$ws.Range("K31").Value = 9451.5
$ws.Range("M31").Value = -9221.5
$ws.Range("H33").Value = 509.84616
$ws.Range("I33").Value = 453.33334
$ws.Range("K33").Value = 453.33334
$ws.Range("M33").Value = -224.33334
$ws.Range("H38").Value = 1575.909
$ws.Range("I38").Value = 259.55554
$ws.Range("J38").Value = 7499.5
$ws.Range("K38").Value = 778.66662
$ws.Range("L38").Value = 22498.5
$ws.Range("M38").Value = -406.66662
$ws.Range("N38").Value = -23242.5
$ws.Range("H51").Value = 12537.375
$ws.Range("J51").Value = 7142.7144
$ws.Range("L51").Value = 7142.7144
$ws.Range("N51").Value = -8110.7144
$ws.Range("H58").Value = 4595
$ws.Range("I58").Value = 1055
$ws.Range("K58").Value = 3165
$ws.Range("M58").Value = -3015
$ws.Range("H106").Value = 7245.875
$ws.Range("I106").Value = 7245.875
$ws.Range("K106").Value = 7245.875
$ws.Range("M106").Value = -6614.875
$ws.Range("H125").Value = 1023.9231
$ws.Range("I125").Value = 339.4
$ws.Range("K125").Value = 3054.6
$ws.Range("M125").Value = -594.5999999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13237.25
$ws.Range("I32").Value = 13523.474
$ws.Range("K32").Value = 13523.474
$ws.Range("M32").Value = -13236.474
$ws.Range("H122").Value = 2070.6667
$ws.Range("I122").Value = 2079.7334
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 6239.2002
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -3789.2002
$ws.Range("N122").Value = -10840
$ws.Range("H132").Value = 6671550.5
$ws.Range("I132").Value = 5387.091
$ws.Range("J132").Value = 25003500
$ws.Range("K132").Value = 16161.273
$ws.Range("L132").Value = 75010500
$ws.Range("M132").Value = -13631.273
$ws.Range("N132").Value = -75015560

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6749.1665
$ws.Range("J86").Value = 9950
$ws.Range("L86").Value = 9950
$ws.Range("N86").Value = -12196
$ws.Range("H89").Value = 6749.1665
$ws.Range("J89").Value = 9950
$ws.Range("L89").Value = 49750
$ws.Range("N89").Value = -60982

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 696.6
$ws.Range("I22").Value = 437.6
$ws.Range("J22").Value = 955.6
$ws.Range("K22").Value = 437.6
$ws.Range("L22").Value = 955.6
$ws.Range("M22").Value = -87.60000000000002
$ws.Range("N22").Value = -1655.6
$ws.Range("H132").Value = 3054.2
$ws.Range("I132").Value = 2504.7778
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 7514.3334
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -4984.3334
$ws.Range("N132").Value = -29057
$ws.Range("H141").Value = 384444.34
$ws.Range("J141").Value = 444285.56
$ws.Range("L141").Value = 444285.56
$ws.Range("N141").Value = -454645.56

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 129.05882
$ws.Range("I2").Value = 119.545456
$ws.Range("J2").Value = 146.5
$ws.Range("K2").Value = 717.272736
$ws.Range("L2").Value = 879
$ws.Range("M2").Value = -604.272736
$ws.Range("N2").Value = -1105
$ws.Range("H41").Value = 12124.8
$ws.Range("I41").Value = 291
$ws.Range("J41").Value = 15083.25
$ws.Range("K41").Value = 873
$ws.Range("L41").Value = 45249.75
$ws.Range("M41").Value = -535
$ws.Range("N41").Value = -45925.75
$ws.Range("H46").Value = 6975.8
$ws.Range("I46").Value = 448.66666
$ws.Range("J46").Value = 16766.5
$ws.Range("K46").Value = 1345.99998
$ws.Range("L46").Value = 50299.5
$ws.Range("M46").Value = -1254.99998
$ws.Range("N46").Value = -50481.5
$ws.Range("H69").Value = 19304.3
$ws.Range("I69").Value = 4099.75
$ws.Range("J69").Value = 29440.666
$ws.Range("K69").Value = 12299.25
$ws.Range("L69").Value = 88321.99800000001
$ws.Range("M69").Value = -11488.25
$ws.Range("N69").Value = -89943.99800000001
$ws.Range("H72").Value = 19304.3
$ws.Range("I72").Value = 4099.75
$ws.Range("J72").Value = 29440.666
$ws.Range("K72").Value = 36897.75
$ws.Range("L72").Value = 264965.994
$ws.Range("M72").Value = -32841.75
$ws.Range("N72").Value = -273077.994
$ws.Range("H76").Value = 20734.5
$ws.Range("I76").Value = 10004
$ws.Range("J76").Value = 25333.285
$ws.Range("K76").Value = 30012
$ws.Range("L76").Value = 75999.855
$ws.Range("M76").Value = -29629
$ws.Range("N76").Value = -76765.855
$ws.Range("H79").Value = 20734.5
$ws.Range("I79").Value = 10004
$ws.Range("J79").Value = 25333.285
$ws.Range("K79").Value = 30012
$ws.Range("L79").Value = 75999.855
$ws.Range("M79").Value = -28686
$ws.Range("N79").Value = -78651.855
$ws.Range("H131").Value = 3882.9429
$ws.Range("J131").Value = 3668.879
$ws.Range("L131").Value = 11006.637
$ws.Range("N131").Value = -21086.637
$ws.Range("H138").Value = 12690.833
$ws.Range("I138").Value = 10112.625
$ws.Range("K138").Value = 30337.875
$ws.Range("M138").Value = -25197.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18079
$ws.Range("J5").Value = 13400
$ws.Range("L5").Value = 13400
$ws.Range("N5").Value = -13624
$ws.Range("H41").Value = 8658
$ws.Range("I41").Value = 3389.6
$ws.Range("K41").Value = 3389.6
$ws.Range("M41").Value = -3034.6
$ws.Range("H80").Value = 5240.4
$ws.Range("I80").Value = 1600
$ws.Range("K80").Value = 1600
$ws.Range("M80").Value = -602
$ws.Range("H83").Value = 5240.4
$ws.Range("I83").Value = 1600
$ws.Range("K83").Value = 8000
$ws.Range("M83").Value = -3008
$ws.Range("H99").Value = 19277.857
$ws.Range("I99").Value = 5824.3335
$ws.Range("K99").Value = 5824.3335
$ws.Range("M99").Value = -3578.3335
$ws.Range("H122").Value = 2808.9443
$ws.Range("I122").Value = 1900.7667
$ws.Range("K122").Value = 5702.300099999999
$ws.Range("M122").Value = -3252.300099999999
$ws.Range("H126").Value = 17163536
$ws.Range("J126").Value = 7570.6665
$ws.Range("L126").Value = 22711.9995
$ws.Range("N126").Value = -27651.9995
$ws.Range("H132").Value = 7074488.5
$ws.Range("I132").Value = 4042.8096
$ws.Range("J132").Value = 31821048
$ws.Range("K132").Value = 12128.4288
$ws.Range("L132").Value = 95463144
$ws.Range("M132").Value = -9598.4288
$ws.Range("N132").Value = -95468204
$ws.Range("H136").Value = 5959.769
$ws.Range("J136").Value = 5959.769
$ws.Range("L136").Value = 17879.307
$ws.Range("N136").Value = -22979.307

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7239.919
$ws.Range("I7").Value = 6720.12
$ws.Range("J7").Value = 8322.833000000001
$ws.Range("K7").Value = 6720.12
$ws.Range("L7").Value = 8322.833000000001
$ws.Range("M7").Value = -6608.12
$ws.Range("N7").Value = -8546.833000000001
$ws.Range("H16").Value = 2345.0908
$ws.Range("I16").Value = 2345.0908
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2345.0908
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2175.0908
$ws.Range("N16").ClearContents()
$ws.Range("H93").Value = 3971368
$ws.Range("I93").Value = 1054.3334
$ws.Range("K93").Value = 1054.3334
$ws.Range("M93").Value = 193.6666
$ws.Range("H126").Value = 7239.919
$ws.Range("I126").Value = 6720.12
$ws.Range("J126").Value = 8322.833000000001
$ws.Range("K126").Value = 20160.36
$ws.Range("L126").Value = 24968.499
$ws.Range("M126").Value = -17690.36
$ws.Range("N126").Value = -29908.499
$ws.Range("H132").Value = 6891.357
$ws.Range("I132").Value = 2888.6
$ws.Range("K132").Value = 8665.799999999999
$ws.Range("M132").Value = -6135.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H100").Value = 1444.5385
$ws.Range("I100").Value = 1453.1111
$ws.Range("K100").Value = 2906.2222
$ws.Range("M100").Value = -2365.2222
$ws.Range("H122").Value = 1807.1515
$ws.Range("I122").Value = 1559.862
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 4679.586
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -2229.586
$ws.Range("N122").Value = -15700
$ws.Range("H132").Value = 557609.0600000001
$ws.Range("I132").Value = 2307.5334
$ws.Range("J132").Value = 3334116.8
$ws.Range("K132").Value = 6922.600199999999
$ws.Range("L132").Value = 10002350.4
$ws.Range("M132").Value = -4392.600199999999
$ws.Range("N132").Value = -10007410.4
$ws.Range("H136").Value = 317092.1
$ws.Range("I136").Value = 4689.185
$ws.Range("J136").Value = 2004067.8
$ws.Range("K136").Value = 14067.555
$ws.Range("L136").Value = 6012203.4
$ws.Range("M136").Value = -11517.555
$ws.Range("N136").Value = -6017303.4
